$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update project title (B3)
$ws.Range("B3").Value = "KLASIFIKASI RONTGEN DADA MENGGUNAKAN DEEP LEARNING"

# Update dataset source link (C3) with hyperlink
$ws.Range("C3").ClearContents()
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.kaggle.com/jtiptj/chest-xray-pneumoniacovid19tuberculosis")

# Update dataset title (D3)
$ws.Range("D3").Value = "Chest X-Ray (Pneumonia, Covid-19, Tuberculosis)"

# Update article source link (E3) with hyperlink
$ws.Range("E3").ClearContents()
$ws.Hyperlinks.Add($ws.Range("E3"), "https://arxiv.org/pdf/2004.05405.pdf")

# Update article title (F3)
$ws.Range("F3").Value = "UNVEILING COVID-19 FROM CHEST X-RAY WITH DEEP LEARNING: A HURDLES RACE WITH SMALL DATA"

$ws.Range("F3").Select()

$wb.Save()
